# CSCockpitBulkTest.xlsx — refresh the bulk return-order regression fixture
# with a new set of order numbers and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testBulkReturnOrders")

# New return-order numbers for rows 2-4 (column A)
$ws.Range("A2").Value = 1000551790
$ws.Range("A3").Value = 1000551792
$ws.Range("A4").Value = 1000551794

# Move/save the active selection on the sheet to C10
$ws.Range("C10").Select()
